$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 133, shifting existing rows 133-161 down to 134-162.
$ws.Rows.Item(133).Insert()

# Populate the newly inserted row 133 with the new record.
$ws.Cells.Item(133, 1).Value = 8
$ws.Cells.Item(133, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(133, 3).Value = "Coquimbo"
$ws.Cells.Item(133, 4).Value = 45077
$ws.Cells.Item(133, 5).Value = 4
$ws.Cells.Item(133, 6).Value = "Fruta"
$ws.Cells.Item(133, 7).Value = 100109
$ws.Cells.Item(133, 8).Value = "Uva"
$ws.Cells.Item(133, 9).Value = 100109001
$ws.Cells.Item(133, 10).Value = "Uva"
$ws.Cells.Item(133, 11).Value = "Autumn Royal"
$ws.Cells.Item(133, 12).Value = "Primera"
$ws.Cells.Item(133, 13).Value = 300
$ws.Cells.Item(133, 14).Value = 13000
$ws.Cells.Item(133, 15).Value = 13500
$ws.Cells.Item(133, 16).Value = 13250
$ws.Cells.Item(133, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(133, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(133, 19).Value = 736
$ws.Cells.Item(133, 20).Value = 18
